$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells stay formatted as text so values
# like "1.001" or "0.000008742" are not reinterpreted as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.037.85'
$ws.Range("E2").Value = '  -2.65%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.822.72'
$ws.Range("E3").Value = '  -1.66%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -1.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.39'
$ws.Range("E5").Value = '  -2.38%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4225'
$ws.Range("E7").Value = '  -1.79%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3674'
$ws.Range("E8").Value = '  -2.26%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07215'
$ws.Range("E9").Value = '  -1.75%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8421'
$ws.Range("E10").Value = '  -4.11%  '
$ws.Range("E11").Value = '  -3.94%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.822.76'
$ws.Range("E12").Value = '  -1.78%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.669'
$ws.Range("E13").Value = '  -1.36%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.07060'
$ws.Range("E14").Value = '  -0.81%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.290'
$ws.Range("E15").Value = '  -2.90%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '89.92'
$ws.Range("E16").Value = '  +1.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.002'
$ws.Range("E17").Value = '  -1.12%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008742'
$ws.Range("E18").Value = '  -2.99%  '
$ws.Range("E19").Value = '  -0.97%  '
$ws.Range("E20").Value = '  -3.78%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.050.56'
$ws.Range("E21").Value = '  -2.70%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.143'
$ws.Range("E22").Value = '  -1.54%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.81'
$ws.Range("E23").Value = '  -2.47%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.045.88'
$ws.Range("E24").Value = '  -1.78%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.980'
$ws.Range("E25").Value = '  -0.69%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.92'
$ws.Range("E26").Value = '  -2.28%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.247'
$ws.Range("E27").Value = '  +3.19%  '
$ws.Range("E28").Value = '  -2.16%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.254'
$ws.Range("E29").Value = '  -2.07%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '116.74'
$ws.Range("E30").Value = '  -1.96%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.175'
$ws.Range("E32").Value = '  -4.51%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7355'
$ws.Range("E33").Value = '  -5.50%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.411'
$ws.Range("E34").Value = '  -3.02%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.879'
$ws.Range("E35").Value = '  -1.83%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.000'
$ws.Range("E36").Value = '  -1.14%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.090'
$ws.Range("E37").Value = '  -3.78%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05259'
$ws.Range("E38").Value = '  -1.69%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01945'
$ws.Range("E39").Value = '  -2.03%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.303'
$ws.Range("E40").Value = '  +1.94%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.867'
$ws.Range("E41").Value = '  -1.27%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1685'
$ws.Range("E42").Value = '  -0.79%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5050'
$ws.Range("E43").Value = '  -1.70%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.542'
$ws.Range("E44").Value = '  -3.40%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.49'
$ws.Range("E45").Value = '  -2.24%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '106.15'
$ws.Range("E46").Value = '  -1.47%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4712'
$ws.Range("E47").Value = '  -1.19%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.937'
$ws.Range("E48").Value = '  +4.97%  '
$ws.Range("E49").Value = '  -1.18%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06321'
$ws.Range("E50").Value = '  -2.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.649'
$ws.Range("E51").Value = '  -2.67%  '
